$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 3014.9167
$ws.Range("I69").Value = 2000
$ws.Range("J69").Value = 3107.182
$ws.Range("K69").Value = 6000
$ws.Range("L69").Value = 9321.545999999998
$ws.Range("M69").Value = -5126
$ws.Range("N69").Value = -11069.546

$ws.Range("H72").Value = 3014.9167
$ws.Range("I72").Value = 2000
$ws.Range("J72").Value = 3107.182
$ws.Range("K72").Value = 18000
$ws.Range("L72").Value = 27964.638
$ws.Range("M72").Value = -13632
$ws.Range("N72").Value = -36700.638

$ws.Range("H76").Value = 85158.63
$ws.Range("I76").Value = 113819.1
$ws.Range("J76").Value = 3271.5715
$ws.Range("K76").Value = 113819.1
$ws.Range("L76").Value = 3271.5715
$ws.Range("M76").Value = -113504.1
$ws.Range("N76").Value = -3901.5715

$ws.Range("H79").Value = 85158.63
$ws.Range("I79").Value = 113819.1
$ws.Range("J79").Value = 3271.5715
$ws.Range("K79").Value = 113819.1
$ws.Range("L79").Value = 3271.5715
$ws.Range("M79").Value = -112727.1
$ws.Range("N79").Value = -5455.5715

$ws.Range("H112").Value = 1447.4814
$ws.Range("J112").Value = 1575.3334
$ws.Range("L112").Value = 4726.0002
$ws.Range("N112").Value = -6942.0002

$ws.Range("H127").Value = 1089
$ws.Range("I127").Value = 878
$ws.Range("J127").Value = 1300
$ws.Range("K127").Value = 2634
$ws.Range("L127").Value = 3900
$ws.Range("M127").Value = 2326
$ws.Range("N127").Value = -13820

$ws.Range("H129").Value = 848.5625
$ws.Range("J129").Value = 1101.2727
$ws.Range("L129").Value = 3303.8181
$ws.Range("N129").Value = -13303.8181

$ws.Range("H137").Value = 23257326
$ws.Range("I137").Value = 1205.1666
$ws.Range("J137").Value = 76925300
$ws.Range("K137").Value = 3615.4998
$ws.Range("L137").Value = 230775900
$ws.Range("M137").Value = -1065.4998
$ws.Range("N137").Value = -230781000

$ws.Range("H138").Value = 1758.3611
$ws.Range("I138").Value = 585.8333
$ws.Range("J138").Value = 3399.9
$ws.Range("K138").Value = 1757.4999
$ws.Range("L138").Value = 10199.7
$ws.Range("M138").Value = 3382.5001
$ws.Range("N138").Value = -20479.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8459.58
$ws.Range("I32").Value = 3215.3513
$ws.Range("J32").Value = 23385.46
$ws.Range("K32").Value = 3215.3513
$ws.Range("L32").Value = 23385.46
$ws.Range("M32").Value = -2928.3513
$ws.Range("N32").Value = -23959.46

$ws.Range("H45").Value = 1987.6
$ws.Range("I45").Value = 2078.5454
$ws.Range("J45").Value = 1737.5
$ws.Range("K45").Value = 2078.5454
$ws.Range("L45").Value = 1737.5
$ws.Range("M45").Value = -1701.5454
$ws.Range("N45").Value = -2491.5

$ws.Range("H74").Value = 5954992
$ws.Range("I74").Value = 9261594
$ws.Range("J74").Value = 3107.5334
$ws.Range("K74").Value = 9261594
$ws.Range("L74").Value = 3107.5334
$ws.Range("M74").Value = -9260720
$ws.Range("N74").Value = -4855.5334

$ws.Range("H77").Value = 5954992
$ws.Range("I77").Value = 9261594
$ws.Range("J77").Value = 3107.5334
$ws.Range("K77").Value = 46307970
$ws.Range("L77").Value = 15537.667
$ws.Range("M77").Value = -46303602
$ws.Range("N77").Value = -24273.667

$ws.Range("H97").Value = 1748.9231
$ws.Range("I97").Value = 2151.5789
$ws.Range("J97").Value = 656
$ws.Range("K97").Value = 2151.5789
$ws.Range("L97").Value = 656
$ws.Range("M97").Value = -1655.5789
$ws.Range("N97").Value = -1648

$ws.Range("H132").Value = 3482.4614
$ws.Range("I132").Value = 3149.1428
$ws.Range("J132").Value = 3871.3333
$ws.Range("K132").Value = 9447.428400000001
$ws.Range("L132").Value = 11613.9999
$ws.Range("M132").Value = -6917.428400000001
$ws.Range("N132").Value = -16673.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2066.6667
$ws.Range("I99").Value = 2900
$ws.Range("J99").Value = 1650
$ws.Range("K99").Value = 2900
$ws.Range("L99").Value = 1650
$ws.Range("M99").Value = -1402
$ws.Range("N99").Value = -4646

$ws.Range("H134").Value = 1171.6346
$ws.Range("I134").Value = 1033.5227
$ws.Range("J134").Value = 1931.25
$ws.Range("K134").Value = 3100.5681
$ws.Range("L134").Value = 5793.75
$ws.Range("M134").Value = -565.5681
$ws.Range("N134").Value = -10863.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3243461.8
$ws.Range("I31").Value = 4370541
$ws.Range("J31").Value = 3108.75
$ws.Range("K31").Value = 4370541
$ws.Range("L31").Value = 3108.75
$ws.Range("M31").Value = -4370246
$ws.Range("N31").Value = -3698.75

$ws.Range("H34").Value = 3243461.8
$ws.Range("I34").Value = 4370541
$ws.Range("J34").Value = 3108.75
$ws.Range("K34").Value = 4370541
$ws.Range("L34").Value = 3108.75
$ws.Range("M34").Value = -4370339
$ws.Range("N34").Value = -3512.75

$ws.Range("H58").Value = 1282.4445
$ws.Range("I58").Value = 685.6316
$ws.Range("J58").Value = 1949.4706
$ws.Range("K58").Value = 685.6316
$ws.Range("L58").Value = 1949.4706
$ws.Range("M58").Value = -482.6316
$ws.Range("N58").Value = -2355.4706

$ws.Range("H132").Value = 1595.75
$ws.Range("I132").Value = 1302.6666
$ws.Range("J132").Value = 2475
$ws.Range("K132").Value = 3907.9998
$ws.Range("L132").Value = 7425
$ws.Range("M132").Value = -1377.9998
$ws.Range("N132").Value = -12485

$ws.Range("H134").Value = 2874.4
$ws.Range("I134").Value = 3237.3333
$ws.Range("J134").Value = 2027.5555
$ws.Range("K134").Value = 9711.999899999999
$ws.Range("L134").Value = 6082.666499999999
$ws.Range("M134").Value = -7176.999899999999
$ws.Range("N134").Value = -11152.6665

$ws.Range("H136").Value = 1282.4445
$ws.Range("I136").Value = 685.6316
$ws.Range("J136").Value = 1949.4706
$ws.Range("K136").Value = 2056.8948
$ws.Range("L136").Value = 5848.4118
$ws.Range("M136").Value = 493.1052
$ws.Range("N136").Value = -10948.4118

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 1214.5714
$ws.Range("I32").Value = 500
$ws.Range("J32").Value = 1333.6666
$ws.Range("K32").Value = 1500
$ws.Range("L32").Value = 4000.9998
$ws.Range("M32").Value = -1217
$ws.Range("N32").Value = -4566.9998

$ws.Range("H35").Value = 3250

$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0

$ws.Range("H131").Value = 1031.7457
$ws.Range("I131").Value = 411.2857
$ws.Range("J131").Value = 1115.2693
$ws.Range("K131").Value = 1233.8571
$ws.Range("L131").Value = 3345.8079
$ws.Range("M131").Value = 3806.1429
$ws.Range("N131").Value = -13425.8079

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2519.0588
$ws.Range("I102").Value = 2579.625
$ws.Range("J102").Value = 1550
$ws.Range("K102").Value = 2579.625
$ws.Range("L102").Value = 1550
$ws.Range("M102").Value = -957.625
$ws.Range("N102").Value = -4794

$ws.Range("H122").Value = 4482.4
$ws.Range("I122").Value = 3637.3333
$ws.Range("J122").Value = 5750
$ws.Range("K122").Value = 10911.9999
$ws.Range("L122").Value = 17250
$ws.Range("M122").Value = -8461.999899999999
$ws.Range("N122").Value = -22150

$ws.Range("H126").Value = 1965.2
$ws.Range("I126").Value = 1859.2632
$ws.Range("J126").Value = 2300.6667
$ws.Range("K126").Value = 5577.7896
$ws.Range("L126").Value = 6902.000100000001
$ws.Range("M126").Value = -3107.7896
$ws.Range("N126").Value = -11842.0001

$ws.Range("H132").Value = 2323.8125
$ws.Range("I132").Value = 2140.0833
$ws.Range("J132").Value = 2875
$ws.Range("K132").Value = 6420.249899999999
$ws.Range("L132").Value = 8625
$ws.Range("M132").Value = -3890.249899999999
$ws.Range("N132").Value = -13685

$ws.Range("H134").Value = 36666
$ws.Range("J134").Value = 36666
$ws.Range("L134").Value = 109998
$ws.Range("N134").Value = -115068

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H96").Value = 26000
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 26000
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 26000
$ws.Range("N96").Value = -31492

$ws.Range("H132").Value = 12506018
$ws.Range("I132").Value = 36778064
$ws.Range("J132").Value = 2236.9092
$ws.Range("K132").Value = 110334192
$ws.Range("L132").Value = 6710.7276
$ws.Range("M132").Value = -110331662
$ws.Range("N132").Value = -11770.7276

$ws.Range("H135").Value = 55360
$ws.Range("J135").Value = 55360
$ws.Range("L135").Value = 55360
$ws.Range("N135").Value = -65500

$ws.Range("H136").Value = 3109.7288
$ws.Range("I136").Value = 3834.6
$ws.Range("J136").Value = 1583.6842
$ws.Range("K136").Value = 11503.8
$ws.Range("L136").Value = 4751.0526
$ws.Range("M136").Value = -8953.799999999999
$ws.Range("N136").Value = -9851.052599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H99").Value = 58333.332
$ws.Range("J99").Value = 58333.332
$ws.Range("L99").Value = 58333.332
$ws.Range("N99").Value = -64323.332

$ws.Range("H107").Value = 909
$ws.Range("I107").Value = 454.44446
$ws.Range("J107").Value = 5000
$ws.Range("K107").Value = 1363.33338
$ws.Range("L107").Value = 15000
$ws.Range("M107").Value = 556.66662
$ws.Range("N107").Value = -18840

$ws.Range("H126").Value = 3004.8262
$ws.Range("I126").Value = 3105.55
$ws.Range("J126").Value = 2333.3333
$ws.Range("K126").Value = 9316.650000000001
$ws.Range("L126").Value = 6999.999899999999
$ws.Range("M126").Value = -6846.650000000001
$ws.Range("N126").Value = -11939.9999
